$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3585.6428
$ws.Range("I40").Value = 6275
$ws.Range("J40").Value = 2509.9
$ws.Range("K40").Value = 6275
$ws.Range("L40").Value = 2509.9
$ws.Range("M40").Value = -6100
$ws.Range("N40").Value = -2859.9
$ws.Range("H43").Value = 3384.111
$ws.Range("I43").Value = 2254.2727
$ws.Range("J43").Value = 5159.5713
$ws.Range("K43").Value = 2254.2727
$ws.Range("L43").Value = 5159.5713
$ws.Range("M43").Value = -2185.2727
$ws.Range("N43").Value = -5297.5713
$ws.Range("H86").Value = 1202
$ws.Range("I86").Value = 1202
$ws.Range("K86").Value = 1202
$ws.Range("M86").Value = -79
$ws.Range("H89").Value = 1202
$ws.Range("I89").Value = 1202
$ws.Range("K89").Value = 6010
$ws.Range("M89").Value = -394
$ws.Range("H96").Value = 52632430
$ws.Range("I96").Value = 58824476
$ws.Range("J96").Value = 73.5
$ws.Range("K96").Value = 176473428
$ws.Range("L96").Value = 220.5
$ws.Range("M96").Value = -176472055
$ws.Range("N96").Value = -2966.5
$ws.Range("H132").Value = 1091.3572
$ws.Range("I132").Value = 1225.909
$ws.Range("J132").Value = 598
$ws.Range("K132").Value = 3677.727
$ws.Range("L132").Value = 1794
$ws.Range("M132").Value = -1147.727
$ws.Range("N132").Value = -6854
$ws.Range("H137").Value = 1948
$ws.Range("I137").Value = 1820.8572
$ws.Range("J137").Value = 2096.3333
$ws.Range("K137").Value = 5462.571599999999
$ws.Range("L137").Value = 6288.999899999999
$ws.Range("M137").Value = -2912.571599999999
$ws.Range("N137").Value = -11388.9999
$ws.Range("H138").Value = 2244.9822
$ws.Range("I138").Value = 1863.0385
$ws.Range("J138").Value = 2576
$ws.Range("K138").Value = 5589.1155
$ws.Range("L138").Value = 7728
$ws.Range("M138").Value = -449.1154999999999
$ws.Range("N138").Value = -18008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 23810700
$ws.Range("I2").Value = 30303584
$ws.Range("J2").Value = 3465
$ws.Range("K2").Value = 30303584
$ws.Range("L2").Value = 3465
$ws.Range("M2").Value = -30303471
$ws.Range("N2").Value = -3691
$ws.Range("H45").Value = 3488.4443
$ws.Range("I45").Value = 2499.25
$ws.Range("J45").Value = 4279.8
$ws.Range("K45").Value = 2499.25
$ws.Range("L45").Value = 4279.8
$ws.Range("M45").Value = -2122.25
$ws.Range("N45").Value = -5033.8
$ws.Range("H56").Value = 27500
$ws.Range("J56").Value = 30000
$ws.Range("L56").Value = 30000
$ws.Range("N56").Value = -31484
$ws.Range("H61").Value = 2018.1305
$ws.Range("I61").Value = 2028.5454
$ws.Range("K61").Value = 2028.5454
$ws.Range("M61").Value = -1816.5454
$ws.Range("H116").Value = 23810700
$ws.Range("I116").Value = 30303584
$ws.Range("J116").Value = 3465
$ws.Range("K116").Value = 30303584
$ws.Range("L116").Value = 3465
$ws.Range("M116").Value = -30301290
$ws.Range("N116").Value = -8053
$ws.Range("H132").Value = 4900
$ws.Range("I132").Value = 4900
$ws.Range("K132").Value = 14700
$ws.Range("M132").Value = -12170
$ws.Range("H136").Value = 2018.1305
$ws.Range("I136").Value = 2028.5454
$ws.Range("K136").Value = 6085.6362
$ws.Range("M136").Value = -3535.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 23810700
$ws.Range("I3").Value = 30303584
$ws.Range("J3").Value = 3465
$ws.Range("K3").Value = 30303584
$ws.Range("L3").Value = 3465
$ws.Range("M3").Value = -30303470
$ws.Range("N3").Value = -3693
$ws.Range("H7").Value = 561.75
$ws.Range("I7").Value = 601
$ws.Range("J7").Value = 522.5
$ws.Range("K7").Value = 601
$ws.Range("L7").Value = 522.5
$ws.Range("M7").Value = -488
$ws.Range("N7").Value = -748.5
$ws.Range("H22").Value = 99.333336
$ws.Range("I22").Value = 99.333336
$ws.Range("K22").Value = 99.333336
$ws.Range("M22").Value = 73.666664
$ws.Range("H94").Value = 2072.3333
$ws.Range("I94").Value = 1454.625
$ws.Range("J94").Value = 3307.75
$ws.Range("K94").Value = 1454.625
$ws.Range("L94").Value = 3307.75
$ws.Range("M94").Value = -1003.625
$ws.Range("N94").Value = -4209.75
$ws.Range("H99").Value = 1313.96
$ws.Range("I99").Value = 1179.0454
$ws.Range("J99").Value = 2303.3333
$ws.Range("K99").Value = 1179.0454
$ws.Range("L99").Value = 2303.3333
$ws.Range("M99").Value = 318.9546
$ws.Range("N99").Value = -5299.3333
$ws.Range("H134").Value = 2447.6667
$ws.Range("I134").Value = 1718.5714
$ws.Range("K134").Value = 5155.7142
$ws.Range("M134").Value = -2620.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3591.6667
$ws.Range("I31").Value = 2400.0527
$ws.Range("K31").Value = 2400.0527
$ws.Range("M31").Value = -2105.0527
$ws.Range("H34").Value = 3591.6667
$ws.Range("I34").Value = 2400.0527
$ws.Range("K34").Value = 2400.0527
$ws.Range("M34").Value = -2198.0527
$ws.Range("H50").Value = 20248.25
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 21712.285
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 21712.285
$ws.Range("M50").Value = -9375
$ws.Range("N50").Value = -22962.285
$ws.Range("H60").Value = 20323
$ws.Range("I60").Value = 10435.556
$ws.Range("J60").Value = 49985.332
$ws.Range("K60").Value = 10435.556
$ws.Range("L60").Value = 49985.332
$ws.Range("M60").Value = -9924.556
$ws.Range("N60").Value = -51007.332
$ws.Range("H132").Value = 1647.375
$ws.Range("I132").Value = 1173.2941
$ws.Range("J132").Value = 2798.7144
$ws.Range("K132").Value = 3519.8823
$ws.Range("L132").Value = 8396.143199999999
$ws.Range("M132").Value = -989.8823000000002
$ws.Range("N132").Value = -13456.1432
$ws.Range("H135").Value = 69
$ws.Range("J135").Value = 69
$ws.Range("L135").Value = 69
$ws.Range("N135").Value = -10209

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1047733
$ws.Range("I4").Value = 1780525
$ws.Range("K4").Value = 5341575
$ws.Range("M4").Value = -5341463
$ws.Range("H7").Value = 50000050
$ws.Range("I7").Value = 50000050
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 150000150
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -150000038
$ws.Range("H37").Value = 98666.664
$ws.Range("J37").Value = 98666.664
$ws.Range("L37").Value = 295999.992
$ws.Range("N37").Value = -296223.992
$ws.Range("H92").Value = 264.44446
$ws.Range("I92").Value = 277.8
$ws.Range("K92").Value = 833.4000000000001
$ws.Range("M92").Value = 414.5999999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6365.5454
$ws.Range("I102").Value = 6669.6665
$ws.Range("J102").Value = 4997
$ws.Range("K102").Value = 6669.6665
$ws.Range("L102").Value = 4997
$ws.Range("M102").Value = -5047.6665
$ws.Range("N102").Value = -8241
$ws.Range("H107").Value = 4190.4
$ws.Range("I107").Value = 2817.3333
$ws.Range("K107").Value = 2817.3333
$ws.Range("M107").Value = -897.3332999999998
$ws.Range("H126").Value = 4933.7
$ws.Range("I126").Value = 4662.4287
$ws.Range("K126").Value = 13987.2861
$ws.Range("M126").Value = -11517.2861
$ws.Range("H132").Value = 2175
$ws.Range("I132").Value = 1810
$ws.Range("K132").Value = 5430
$ws.Range("M132").Value = -2900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4746.7
$ws.Range("I22").Value = 6168.5
$ws.Range("K22").Value = 6168.5
$ws.Range("M22").Value = -5873.5
$ws.Range("H27").Value = 4746.7
$ws.Range("I27").Value = 6168.5
$ws.Range("K27").Value = 6168.5
$ws.Range("M27").Value = -6061.5
$ws.Range("H46").Value = 252994.25
$ws.Range("I46").Value = 3992.6667
$ws.Range("J46").Value = 999999
$ws.Range("K46").Value = 3992.6667
$ws.Range("L46").Value = 999999
$ws.Range("M46").Value = -3804.6667
$ws.Range("N46").Value = -1000375
$ws.Range("H55").Value = 873.625
$ws.Range("J55").Value = 1199
$ws.Range("L55").Value = 1199
$ws.Range("N55").Value = -1545
$ws.Range("H82").Value = 906.5294
$ws.Range("I82").Value = 616.46155
$ws.Range("J82").Value = 1849.25
$ws.Range("K82").Value = 616.46155
$ws.Range("L82").Value = 1849.25
$ws.Range("M82").Value = -255.46155
$ws.Range("N82").Value = -2571.25
$ws.Range("H85").Value = 906.5294
$ws.Range("I85").Value = 616.46155
$ws.Range("J85").Value = 1849.25
$ws.Range("K85").Value = 616.46155
$ws.Range("L85").Value = 1849.25
$ws.Range("M85").Value = 631.53845
$ws.Range("N85").Value = -4345.25
$ws.Range("H132").Value = 4464.6665
$ws.Range("J132").Value = 4444
$ws.Range("L132").Value = 13332
$ws.Range("N132").Value = -18392

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 16739.5
$ws.Range("J54").Value = 17296
$ws.Range("L54").Value = 17296
$ws.Range("N54").Value = -18336
$ws.Range("H81").Value = 3438.6667
$ws.Range("I81").Value = 3493.5
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 6987
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -5926
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 3438.6667
$ws.Range("I84").Value = 3493.5
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 34935
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -29631
$ws.Range("N84").Value = -40608
$ws.Range("H107").Value = 494
$ws.Range("I107").Value = 491.75
$ws.Range("J107").Value = 503
$ws.Range("K107").Value = 1475.25
$ws.Range("L107").Value = 1509
$ws.Range("M107").Value = 444.75
$ws.Range("N107").Value = -5349
$ws.Range("H132").Value = 2525.25
$ws.Range("I132").Value = 1813
$ws.Range("K132").Value = 5439
$ws.Range("M132").Value = -2909
